# Change all units to SI/pint-compatible: "EJ/y" -> "EJ"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "EJ"
$ws.Range("D3").Value = "EJ"

# Update the view state: select D4 and scroll so A4 is the top-left visible cell
$ws.Activate()
$ws.Range("D4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
